# Qatar Stars League workbook update, 26-02-2024 22:04
# The source feed re-ordered several fixtures that share the same kickoff
# date/time; this swaps the full data row (columns B:AC, i.e. everything
# except the running index in column A) between each affected pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(27, 28),
    @(40, 41),
    @(42, 43),
    @(45, 46),
    @(56, 57),
    @(75, 76),
    @(78, 79)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Write-Output "Row pairs swapped."
